# Add the "currency" ("Moeda") field to the "contract" sheet.
# This introduces a new column N with a header and one value per
# existing data row (rows 2-5), mirroring the shared-strings / sheet1
# changes from the commit "Including the field "currency" on model
# "contract"".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "Moeda"
$ws.Range("N2").Value = "Real"
$ws.Range("N3").Value = "Real"
$ws.Range("N4").Value = "Dolar"
$ws.Range("N5").Value = "Euro"

# Match the author's final selection/view state as closely as possible.
$ws.Range("N4").Select()
